$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel constant used below:
#   xlShiftToRight = -4161  (used when inserting a copied column)
$xlShiftToRight = -4161

# The table currently ends at column Q (year 2020). Two more year columns
# (2021, 2022) need to be appended. The cleanest way to make the new
# columns R and S inherit exactly the same cell formatting/styles as the
# existing 2020 column Q is to copy column Q and insert the copy at R,
# then again at S - this duplicates both values and formatting, after
# which we overwrite the values with the real 2021/2022 figures.
$ws.Columns("Q").Copy() | Out-Null
$ws.Columns("R").Insert($xlShiftToRight) | Out-Null
$ws.Columns("Q").Copy() | Out-Null
$ws.Columns("S").Insert($xlShiftToRight) | Out-Null
$ws.Application.CutCopyMode = $false

# --- Row 4 (years header): 2021 / 2022 ---
$ws.Range("R4").Value = 2021
$ws.Range("S4").Value = 2022

# --- Row 5 (Доходы, всего / Revenues, total) ---
$ws.Range("P5").Value = 25.6
$ws.Range("Q5").Value = 23.8
$ws.Range("R5").Value = 26.8
$ws.Range("S5").Value = 26.8

# --- Row 6 (Налоговые доходы / Tax revenues) ---
$ws.Range("P6").Value = 18.600000000000001
$ws.Range("Q6").Value = 16.7
$ws.Range("R6").Value = 19.3
$ws.Range("S6").Value = 19.3

# --- Row 7 (Received official transfers) ---
$ws.Range("R7").Value = "-"
$ws.Range("S7").Value = "-"

# --- Row 8 (Неналоговые доходы / Non-tax revenues) ---
$ws.Range("P8").Value = 2.1
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 1.8
$ws.Range("S8").Value = 1.8

# --- Row 9 (Revenues from the sale of non-financial assets) ---
$ws.Range("P9").Value = 4.9000000000000004
$ws.Range("Q9").Value = 5.2
$ws.Range("R9").Value = 5.7
$ws.Range("S9").Value = 5.7

# --- Row 10 (last data row) ---
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0

# Selection moves to T3 in the saved file
$ws.Range("T3").Select() | Out-Null
